# (#33) Alteração nos rótulos da tabela para já transformar a primeira linha
# em cabeçalho automaticamente no Power BI.
#
# Prefixa os rótulos de ano/intervalo da primeira linha (cabeçalho) de cada
# planilha com "Ano " ou "Intervalo ", conforme o caso.

$wb = $excel.ActiveWorkbook

# Planilhas cujo cabeçalho usa o prefixo "Ano " nas colunas B:E (ou só B).
$anoSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)",
    "Custo Total (bilhões de R$)"
)

foreach ($sheetName in $anoSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($col = 2; $col -le 5; $col++) {
        $cell = $ws.Cells.Item(1, $col)
        $val = $cell.Value2
        if ($val -ne $null -and $val -ne "") {
            $cell.Value = "Ano " + $val
        }
    }
}

# Planilha cujo cabeçalho usa o prefixo "Intervalo ".
$ws4 = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
for ($col = 2; $col -le 5; $col++) {
    $cell = $ws4.Cells.Item(1, $col)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $cell.Value = "Intervalo " + $val
    }
}
